# Ccl12-Ackr4.xlsx: refresh NATMI TPM numbers and cell-type pairing.
# - New "MuSCs" sending-cluster rows appear; the old Inflammatory-Mac/
#   Neutrophils/Resolving-Mac x Neutrophils target-cluster row is gone
#   (9 data rows -> 8 data rows), so the sheet shrinks from A1:T10 to A1:T9.
# - Every numeric column (E:T) is recomputed with the new TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old last row (row 10) no longer exists in the refreshed data - remove it.
$ws.Rows(10).Delete()

# Row 2: Inflammatory-Mac -> Ccl12/Ackr4 -> ECs
$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("B2").Value = "Ccl12"
$ws.Range("C2").Value = "Ackr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 42.70362466666668
$ws.Range("H2").Value = 128.110874
$ws.Range("I2").Value = 0.510021191154308
$ws.Range("J2").Value = 0.5102913077099245
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.1631145
$ws.Range("N2").Value = 0.326229
$ws.Range("O2").Value = 0.7212828052797984
$ws.Range("P2").Value = 0.7212828052797984
$ws.Range("Q2").Value = 6.965580385691002
$ws.Range("R2").Value = 41.793482314146
$ws.Range("S2").Value = 0.3678695155079236
$ws.Range("T2").Value = 0.3680643459349111

# Row 3: Inflammatory-Mac -> Ccl12/Ackr4 -> MuSCs
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Ccl12"
$ws.Range("C3").Value = "Ackr4"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 42.70362466666668
$ws.Range("H3").Value = 128.110874
$ws.Range("I3").Value = 0.510021191154308
$ws.Range("J3").Value = 0.5102913077099245
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.0630305
$ws.Range("N3").Value = 0.126061
$ws.Range("O3").Value = 0.2787171947202017
$ws.Range("P3").Value = 0.2787171947202017
$ws.Range("Q3").Value = 2.691630814552334
$ws.Range("R3").Value = 16.149784887314
$ws.Range("S3").Value = 0.1421516756463845
$ws.Range("T3").Value = 0.1422269617750134

# Row 4: MuSCs -> Ccl12/Ackr4 -> ECs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Ccl12"
$ws.Range("C4").Value = "Ackr4"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.132963
$ws.Range("H4").Value = 0.265926
$ws.Range("I4").Value = 0.001588013855235666
$ws.Range("J4").Value = 0.001059236597621443
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.1631145
$ws.Range("N4").Value = 0.326229
$ws.Range("O4").Value = 0.7212828052797984
$ws.Range("P4").Value = 0.7212828052797984
$ws.Range("Q4").Value = 0.0216881932635
$ws.Range("R4").Value = 0.08675277305399999
$ws.Range("S4").Value = 0.001145407088327569
$ws.Range("T4").Value = 0.0007640091445874233

# Row 5: MuSCs -> Ccl12/Ackr4 -> MuSCs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Ccl12"
$ws.Range("C5").Value = "Ackr4"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.132963
$ws.Range("H5").Value = 0.265926
$ws.Range("I5").Value = 0.001588013855235666
$ws.Range("J5").Value = 0.001059236597621443
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = 0.0630305
$ws.Range("N5").Value = 0.126061
$ws.Range("O5").Value = 0.2787171947202017
$ws.Range("P5").Value = 0.2787171947202017
$ws.Range("Q5").Value = 0.0083807243715
$ws.Range("R5").Value = 0.033522897486
$ws.Range("S5").Value = 0.0004426067669080974
$ws.Range("T5").Value = 0.0002952274530340196

# Row 6: Neutrophils -> Ccl12/Ackr4 -> ECs
$ws.Range("A6").Value = "Neutrophils"
$ws.Range("B6").Value = "Ccl12"
$ws.Range("C6").Value = "Ackr4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 8.925701999999999
$ws.Range("H6").Value = 26.777106
$ws.Range("I6").Value = 0.1066021257320059
$ws.Range("J6").Value = 0.1066585841684857
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.5
$ws.Range("M6").Value = 0.1631145
$ws.Range("N6").Value = 0.326229
$ws.Range("O6").Value = 0.7212828052797984
$ws.Range("P6").Value = 0.7212828052797984
$ws.Range("Q6").Value = 1.455911418879
$ws.Range("R6").Value = 8.735468513274
$ws.Range("S6").Value = 0.07689028029677097
$ws.Range("T6").Value = 0.07693100279621684

# Row 7: Neutrophils -> Ccl12/Ackr4 -> MuSCs
$ws.Range("A7").Value = "Neutrophils"
$ws.Range("B7").Value = "Ccl12"
$ws.Range("C7").Value = "Ackr4"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 8.925701999999999
$ws.Range("H7").Value = 26.777106
$ws.Range("I7").Value = 0.1066021257320059
$ws.Range("J7").Value = 0.1066585841684857
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.5
$ws.Range("M7").Value = 0.0630305
$ws.Range("N7").Value = 0.126061
$ws.Range("O7").Value = 0.2787171947202017
$ws.Range("P7").Value = 0.2787171947202017
$ws.Range("Q7").Value = 0.562591459911
$ws.Range("R7").Value = 3.375548759466
$ws.Range("S7").Value = 0.0297118454352349
$ws.Range("T7").Value = 0.02972758137226884

# Row 8: Resolving-Mac -> Ccl12/Ackr4 -> ECs
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Ccl12"
$ws.Range("C8").Value = "Ackr4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 31.96682866666667
$ws.Range("H8").Value = 95.900486
$ws.Range("I8").Value = 0.3817886692584505
$ws.Range("J8").Value = 0.3819908715239683
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.5
$ws.Range("M8").Value = 0.1631145
$ws.Range("N8").Value = 0.326229
$ws.Range("O8").Value = 0.7212828052797984
$ws.Range("P8").Value = 0.7212828052797984
$ws.Range("Q8").Value = 5.214253274549
$ws.Range("R8").Value = 31.285519647294
$ws.Range("S8").Value = 0.2753776023867763
$ws.Range("T8").Value = 0.2755234474040829

# Row 9: Resolving-Mac -> Ccl12/Ackr4 -> MuSCs
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Ccl12"
$ws.Range("C9").Value = "Ackr4"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 31.96682866666667
$ws.Range("H9").Value = 95.900486
$ws.Range("I9").Value = 0.3817886692584505
$ws.Range("J9").Value = 0.3819908715239683
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.5
$ws.Range("M9").Value = 0.0630305
$ws.Range("N9").Value = 0.126061
$ws.Range("O9").Value = 0.2787171947202017
$ws.Range("P9").Value = 0.2787171947202017
$ws.Range("Q9").Value = 2.014885194274334
$ws.Range("R9").Value = 12.089311165646
$ws.Range("S9").Value = 0.1064110668716742
$ws.Range("T9").Value = 0.1064674241198854
